# Version 0.2.10: Bug fix #14. New plugin: SUEWSAnalyzer. General Bug fixing
#
# Applies the recorded edit to SUEWS_init.xlsx:
#  - fixes the "Deciduous" sheet/label spelling to "Decidious"
#  - bumps the per-sheet SiteList row-count index in column F
#  - resets each sheet's remembered cell selection
#  - moves the active tab / tabSelected flag to "ESTM Coefficients"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheet 5 "Deciduous" -> "Decidious" (typo fix) and correct the
#    surface-characteristics label cell that lives on that sheet.
# ---------------------------------------------------------------------------
$wsDecidious = $wb.Worksheets.Item("Deciduous")
$wsDecidious.Name = "Decidious"
$wsDecidious.Range("C1").Value = "Decidious surface characteristics"
$wsDecidious.Range("F1").Value = 37
$wsDecidious.Range("F2").Select()

# ---------------------------------------------------------------------------
# 2) Per-sheet column F counter bumps + selection resets.
#    (sheet1 "Example DO NOT CHANGE POSITION" is untouched by the edit.)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Paved")
$ws.Range("F1").Value = 34
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Building")
$ws.Range("F1").Value = 35
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Evergreen")
$ws.Range("F1").Value = 36
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Grass")
$ws.Range("F1").Value = 38
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Bare Soil")
$ws.Range("F1").Value = 39
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Water")
$ws.Range("F1").Value = 40
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Conductance")
$ws.Range("F1").Value = 45
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Snow")
$ws.Range("F1").Value = 46
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Snow clearing")
$ws.Range("F1").Value = 47
$ws.Range("F2").Value = 48
$ws.Range("F3").Select()

$ws = $wb.Worksheets.Item("Anthropogenic")
$ws.Range("F1").Value = 49
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Energy")
$ws.Range("F1").Value = 50
$ws.Range("F2").Value = 51
$ws.Range("F3").Select()

$ws = $wb.Worksheets.Item("Irrigation")
$ws.Range("F1").Value = 54
$ws.Range("F2").Select()

$ws = $wb.Worksheets.Item("Water Use (Manual)")
$ws.Range("F1").Value = 55
$ws.Range("F2").Value = 56
$ws.Range("F3").Select()

$ws = $wb.Worksheets.Item("Water Use (Automatic)")
$ws.Range("F1").Value = 57
$ws.Range("F2").Value = 58
$ws.Range("F3").Select()

$ws = $wb.Worksheets.Item("ESTM Coefficients")
$ws.Range("F1").Value = 85
$ws.Range("F2").Select()

# ---------------------------------------------------------------------------
# 3) Make "ESTM Coefficients" the active tab (this also flips the
#    tabSelected flag from "Bare Soil" onto this sheet) and set the window
#    tab-strip ratio recorded in the source edit.
# ---------------------------------------------------------------------------
$ws.Activate()

$win = $excel.ActiveWindow
$win.TabRatio = 867
